# Quarterly indexing esoteric bug-fix operation
#
# The date stamps in column A (rows 2-73) were computed with an off-by-one
# quarter indexing bug: each date should represent the 15th of the month
# following the recorded quarter-start month, but was instead left at the
# 1st of that quarter-start month. This corrects every date in A2:A73 by
# advancing it one month and moving it to the 15th of that month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date epoch (serial day 0) is 1899-12-30.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$firstRow = 2
$lastRow = 73

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2
    $oldDate = $epoch.AddDays($oldSerial)
    $newDate = $oldDate.AddMonths(1).AddDays(14)
    $newSerial = $newDate.ToOADate()
    $cell.Value = $newSerial
}
